$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (each row 2-10 shifted up from the row below it,
# row 11 receives newly computed ifoCAST full-series evaluation values)
$data = @(
    @(-0.2196048395615291, 2.254442816470436, 18.45943047502218, 4.296443933652827, 4.39180238609223, 22),
    @(-0.6722431955569657, 1.734626294717863, 10.27176916279258, 3.204960087550636, 3.211051410560245, 21),
    @(-0.2619961853625158, 1.574840796137216, 11.52157846430296, 3.394345071483298, 3.472135163382664, 20),
    @(-0.2304210235372018, 1.666997890576107, 10.06653556874298, 3.172780416092955, 3.251114276439574, 19),
    @(-0.2737229439002919, 1.775059182383401, 10.57130412137001, 3.251354198079626, 3.333738739228395, 18),
    @(-0.1641233777288165, 1.781553791812374, 11.51292075854603, 3.393069518672736, 3.493402108638973, 17),
    @(-0.1510345969195566, 1.85455728114967, 11.36885938571983, 3.37177392268815, 3.478857733797603, 16),
    @(-0.1166201009408896, 1.99642669408382, 12.99001923899349, 3.604166927182131, 3.728713727525796, 15),
    @(-0.07784734317373232, 1.990677678875613, 13.16673245956295, 3.628599242071649, 3.764708707805368, 14),
    @(-0.04637512304105267, 2.085006870274704, 13.21266471962553, 3.634922931731226, 3.783039814737283, 13)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
